$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking Price cells so they are not
# auto-converted to Number type (preserves literal string, e.g. trailing zeros).
$textCells = @("D5", "D6", "D9", "D11", "D12", "D15", "D16", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D29", "D30", "D32", "D33", "D34", "D36", "D37", "D38", "D39", "D40", "D44", "D45", "D46", "D47", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '64.451.68'
$ws.Range("E2").Value = '  +2.16%  '
$ws.Range("D3").Value = '3.449.53'
$ws.Range("E3").Value = '  +2.32%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '574.94'
$ws.Range("E5").Value = '  +0.43%  '
$ws.Range("D6").Value = '157.98'
$ws.Range("E6").Value = '  +3.48%  '
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("D8").Value = '3.451.24'
$ws.Range("E8").Value = '  +2.06%  '
$ws.Range("D9").Value = '0.582'
$ws.Range("E9").Value = '  +10.37%  '
$ws.Range("E10").Value = '  -0.64%  '
$ws.Range("D11").Value = '0.125'
$ws.Range("E11").Value = '  +4.84%  '
$ws.Range("D12").Value = '0.445'
$ws.Range("E12").Value = '  +1.99%  '
$ws.Range("D13").Value = '4.035.96'
$ws.Range("E13").Value = '  +2.04%  '
$ws.Range("E14").Value = '  -2.91%  '
$ws.Range("D15").Value = '0.0000194'
$ws.Range("E15").Value = '  +7.71%  '
$ws.Range("D16").Value = '28.22'
$ws.Range("E16").Value = '  +4.52%  '
$ws.Range("D17").Value = '64.416.73'
$ws.Range("E17").Value = '  +2.00%  '
$ws.Range("D18").Value = '3.435.92'
$ws.Range("E18").Value = '  +2.34%  '
$ws.Range("D19").Value = '6.46'
$ws.Range("E19").Value = '  +2.20%  '
$ws.Range("D20").Value = '14.38'
$ws.Range("E20").Value = '  +3.28%  '
$ws.Range("D21").Value = '389.78'
$ws.Range("E21").Value = '  +1.11%  '
$ws.Range("D22").Value = '8.25'
$ws.Range("E22").Value = '  -2.01%  '
$ws.Range("D23").Value = '73.90'
$ws.Range("E23").Value = '  +4.78%  '
$ws.Range("D24").Value = '0.543'
$ws.Range("E24").Value = '  +1.78%  '
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  +0.12%  '
$ws.Range("E26").Value = '  +24.97%  '
$ws.Range("D27").Value = '9.59'
$ws.Range("E27").Value = '  +3.89%  '
$ws.Range("E28").Value = '  +0.36%  '
$ws.Range("D29").Value = '1.01'
$ws.Range("E29").Value = '  +0.66%  '
$ws.Range("D30").Value = '6.22'
$ws.Range("E30").Value = '  +11.79%  '
$ws.Range("E31").Value = '  +10.37%  '
$ws.Range("D32").Value = '2.03'
$ws.Range("E32").Value = '  +0.69%  '
$ws.Range("D33").Value = '6.59'
$ws.Range("D34").Value = '23.68'
$ws.Range("E34").Value = '  +2.59%  '
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("D36").Value = '6.98'
$ws.Range("E36").Value = '  +3.91%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = '1.48'
$ws.Range("E37").Value = '  +0.48%  '
$ws.Range("B38").Value = 'Monero'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D38").Value = '160.69'
$ws.Range("E38").Value = '  +1.48%  '
$ws.Range("D39").Value = '1.90'
$ws.Range("E39").Value = '  +1.74%  '
$ws.Range("D40").Value = '0.0773'
$ws.Range("E40").Value = '  +4.06%  '
$ws.Range("D41").Value = '2.935.64'
$ws.Range("E41").Value = '  +1.78%  '
$ws.Range("E42").Value = '  -0.96%  '
$ws.Range("E43").Value = '  -3.35%  '
$ws.Range("D44").Value = '42.65'
$ws.Range("E44").Value = '  +4.50%  '
$ws.Range("D45").Value = '4.43'
$ws.Range("E45").Value = '  +4.21%  '
$ws.Range("D46").Value = '0.767'
$ws.Range("E46").Value = '  +2.48%  '
$ws.Range("D47").Value = '23.66'
$ws.Range("E47").Value = '  +7.93%  '
$ws.Range("E48").Value = '  +3.82%  '
$ws.Range("D49").Value = '2.24'
$ws.Range("E49").Value = '  +22.82%  '
$ws.Range("E50").Value = '  +4.55%  '
$ws.Range("D51").Value = '0.861'
$ws.Range("E51").Value = '  +6.92%  '
